# "Add questions algorithm corrected"
#
# Before: a single sheet "Sheet" holding the full C++ question bank.
# After : two sheets -
#   1) "C_Bank"    - new, small "C" question bank (header + 1 sample row)
#   2) "C++_Bank"  - the original question bank (renamed), with a handful
#                    of row-height tweaks and the autofilter / hidden
#                    _FilterDatabase name now living on it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the existing sheet so we have one copy to trim down into
#    the new "C_Bank" sheet and one copy that keeps the full original
#    C++ question bank ("C++_Bank"). Copying (rather than building from
#    scratch) keeps every existing style/format index intact.
# ---------------------------------------------------------------------
$original = $wb.Worksheets.Item(1)
$original.Copy($null, $original)

$cBank   = $wb.Worksheets.Item(1)
$cppBank = $wb.Worksheets.Item(2)

$cBank.Name   = "C_Bank"
$cppBank.Name = "C++_Bank"

# ---------------------------------------------------------------------
# 2. Trim "C_Bank" down to the header row plus a single sample question.
# ---------------------------------------------------------------------
$cBank.Rows("3:49").Delete()

# Reuse the existing plain-data style (already applied to C2:E2) for the
# whole data row instead of letting Excel mint a brand-new style index.
$cBank.Range("C2").Copy()
$cBank.Range("A2:F2").PasteSpecial(-4122)   # xlPasteFormats

$cBank.Range("A2").Value = 2
$cBank.Range("B2").Value = "This is sample C question ';"
$cBank.Range("C2").Value = "High"
$cBank.Range("D2:F2").ClearContents()

# Drop the leftover custom row height from the old row 2 (409.5) so it
# falls back to the sheet's default height.
$cBank.Rows(2).AutoFit()

# This sheet never had the C++ bank's filter turned on.
$cBank.AutoFilterMode = $false

# ---------------------------------------------------------------------
# 3. Re-home the (hidden) _FilterDatabase defined name onto "C++_Bank",
#    since that's the sheet that still carries the autoFilter.
# ---------------------------------------------------------------------
$oldFilterName = $wb.Names.Item(1)
$oldFilterName.Delete()
$filterName = $cppBank.Names.Add("_xlnm._FilterDatabase", "='C++_Bank'!`$A`$1:`$E`$8")
$filterName.Visible = $false

# ---------------------------------------------------------------------
# 4. Small row-height rounding corrections on "C++_Bank" (cosmetic
#    re-measurement noise from the original edit, reproduced exactly).
# ---------------------------------------------------------------------
$cppBank.Rows(5).RowHeight  = 16.9
$cppBank.Rows(6).RowHeight  = 403.15
$cppBank.Rows(17).RowHeight = 403.15
$cppBank.Rows(24).RowHeight = 25.9
$cppBank.Rows(29).RowHeight = 100.9
$cppBank.Rows(31).RowHeight = 100.9
$cppBank.Rows(33).RowHeight = 319.15
$cppBank.Rows(38).RowHeight = 268.9
$cppBank.Rows(39).RowHeight = 16.9
$cppBank.Rows(41).RowHeight = 37.9
$cppBank.Rows(44).RowHeight = 16.5
$cppBank.Rows(45).RowHeight = 16.5
$cppBank.Rows(46).RowHeight = 16.5
$cppBank.Rows(47).RowHeight = 16.5
$cppBank.Rows(48).RowHeight = 16.5
$cppBank.Rows(49).RowHeight = 16.5

# ---------------------------------------------------------------------
# 5. Selection / active-sheet bookkeeping to match the saved view state.
# ---------------------------------------------------------------------
$cppBank.Activate()
$cppBank.Range("A1:XFD1").Select()

$cBank.Activate()
$cBank.Range("E14").Select()
